$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 560980.5600000001
$ws.Range("J17").Value = 560980.5600000001
$ws.Range("L17").Value = 1682941.68
$ws.Range("N17").Value = -1683277.68

$ws.Range("H62").Value = 1632.3182
$ws.Range("I62").Value = 1569
$ws.Range("J62").Value = 1650.9412
$ws.Range("K62").Value = 1569
$ws.Range("L62").Value = 1650.9412
$ws.Range("M62").Value = -945
$ws.Range("N62").Value = -2898.9412

$ws.Range("H65").Value = 1632.3182
$ws.Range("I65").Value = 1569
$ws.Range("J65").Value = 1650.9412
$ws.Range("K65").Value = 7845
$ws.Range("L65").Value = 8254.706
$ws.Range("M65").Value = -4725
$ws.Range("N65").Value = -14494.706

$ws.Range("H92").Value = 275.8
$ws.Range("I92").Value = 238.58824
$ws.Range("J92").Value = 486.66666
$ws.Range("K92").Value = 238.58824
$ws.Range("L92").Value = 486.66666
$ws.Range("M92").Value = 1009.41176
$ws.Range("N92").Value = -2982.66666

$ws.Range("H137").Value = 3428.547
$ws.Range("I137").Value = 1046.5
$ws.Range("J137").Value = 6096.44
$ws.Range("K137").Value = 3139.5
$ws.Range("L137").Value = 18289.32
$ws.Range("M137").Value = -589.5
$ws.Range("N137").Value = -23389.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1229.4906
$ws.Range("I61").Value = 777.13336
$ws.Range("J61").Value = 1819.5217
$ws.Range("K61").Value = 777.13336
$ws.Range("L61").Value = 1819.5217
$ws.Range("M61").Value = -565.13336
$ws.Range("N61").Value = -2243.5217

$ws.Range("H74").Value = 2381.6462
$ws.Range("I74").Value = 797.5893
$ws.Range("K74").Value = 797.5893
$ws.Range("M74").Value = 76.41070000000002

$ws.Range("H77").Value = 2381.6462
$ws.Range("I77").Value = 797.5893
$ws.Range("K77").Value = 3987.9465
$ws.Range("M77").Value = 380.0535

$ws.Range("H132").Value = 12554.526
$ws.Range("I132").Value = 9916.154
$ws.Range("K132").Value = 29748.462
$ws.Range("M132").Value = -27218.462

$ws.Range("H136").Value = 1229.4906
$ws.Range("I136").Value = 777.13336
$ws.Range("J136").Value = 1819.5217
$ws.Range("K136").Value = 2331.40008
$ws.Range("L136").Value = 5458.5651
$ws.Range("M136").Value = 218.5999199999997
$ws.Range("N136").Value = -10558.5651

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2101.9375
$ws.Range("I105").Value = 1328
$ws.Range("J105").Value = 2453.7273
$ws.Range("K105").Value = 1328
$ws.Range("L105").Value = 2453.7273
$ws.Range("M105").Value = 419
$ws.Range("N105").Value = -5947.7273

$ws.Range("H134").Value = 1323.1428
$ws.Range("I134").Value = 1159.875
$ws.Range("J134").Value = 1845.6
$ws.Range("K134").Value = 3479.625
$ws.Range("L134").Value = 5536.799999999999
$ws.Range("M134").Value = -944.625
$ws.Range("N134").Value = -10606.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 700
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0

$ws.Range("H31").Value = 36705.15
$ws.Range("I31").Value = 1075.2354
$ws.Range("J31").Value = 81572.44500000001
$ws.Range("K31").Value = 1075.2354
$ws.Range("L31").Value = 81572.44500000001
$ws.Range("M31").Value = -780.2354
$ws.Range("N31").Value = -82162.44500000001

$ws.Range("H34").Value = 36705.15
$ws.Range("I34").Value = 1075.2354
$ws.Range("J34").Value = 81572.44500000001
$ws.Range("K34").Value = 1075.2354
$ws.Range("L34").Value = 81572.44500000001
$ws.Range("M34").Value = -873.2354
$ws.Range("N34").Value = -81976.44500000001

$ws.Range("H99").Value = 2268.182
$ws.Range("I99").Value = 1387.2
$ws.Range("J99").Value = 3002.3333
$ws.Range("K99").Value = 1387.2
$ws.Range("L99").Value = 3002.3333
$ws.Range("M99").Value = 110.8
$ws.Range("N99").Value = -5998.3333

$ws.Range("H126").Value = 2268.182
$ws.Range("I126").Value = 1387.2
$ws.Range("J126").Value = 3002.3333
$ws.Range("K126").Value = 4161.6
$ws.Range("L126").Value = 9006.999899999999
$ws.Range("M126").Value = -1691.6
$ws.Range("N126").Value = -13946.9999

$ws.Range("H132").Value = 18871448
$ws.Range("I132").Value = 23813510
$ws.Range("J132").Value = 1753.4546
$ws.Range("K132").Value = 71440530
$ws.Range("L132").Value = 5260.3638
$ws.Range("M132").Value = -71438000
$ws.Range("N132").Value = -10320.3638

$ws.Range("H134").Value = 5051.1514
$ws.Range("I134").Value = 5989.84
$ws.Range("J134").Value = 2117.75
$ws.Range("K134").Value = 17969.52
$ws.Range("L134").Value = 6353.25
$ws.Range("M134").Value = -15434.52
$ws.Range("N134").Value = -11423.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1586.9375
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 6000
$ws.Range("N4").Value = -6224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4393.2383
$ws.Range("I132").Value = 5508.8623
$ws.Range("J132").Value = 1904.5385
$ws.Range("K132").Value = 16526.5869
$ws.Range("L132").Value = 5713.6155
$ws.Range("M132").Value = -13996.5869
$ws.Range("N132").Value = -10773.6155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 284054.5
$ws.Range("I2").Value = 5000001
$ws.Range("J2").Value = 6645.8823
$ws.Range("K2").Value = 5000001
$ws.Range("L2").Value = 6645.8823
$ws.Range("N2").Value = -6869.8823
$ws.Range("M2").Value = -4999889

$ws.Range("H40").Value = 2122.75
$ws.Range("I40").Value = 1829.7142
$ws.Range("J40").Value = 2533
$ws.Range("K40").Value = 1829.7142
$ws.Range("L40").Value = 2533
$ws.Range("M40").Value = -1693.7142
$ws.Range("N40").Value = -2805

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1000000
$ws.Range("J2").Value = 1000000
$ws.Range("L2").Value = 1000000
$ws.Range("N2").Value = -1000224

$ws.Range("H126").Value = 859.3333
$ws.Range("I126").Value = 739.2
$ws.Range("J126").Value = 1460
$ws.Range("K126").Value = 2217.6
$ws.Range("L126").Value = 4380
$ws.Range("M126").Value = 252.3999999999996
$ws.Range("N126").Value = -9320

$ws.Range("H132").Value = 2138.7974
$ws.Range("I132").Value = 2936.2126
$ws.Range("J132").Value = 967.59375
$ws.Range("K132").Value = 8808.6378
$ws.Range("L132").Value = 2902.78125
$ws.Range("M132").Value = -6278.6378
$ws.Range("N132").Value = -7962.78125

$ws.Range("H136").Value = 1246.3286
$ws.Range("I136").Value = 748.975
$ws.Range("J136").Value = 1909.4667
$ws.Range("K136").Value = 2246.925
$ws.Range("L136").Value = 5728.4001
$ws.Range("M136").Value = 303.0749999999998
$ws.Range("N136").Value = -10828.4001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N4").ClearContents()
